# Applies the cryptos list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.309.80"
Set-TextValue "E2" "  -3.01%  "

Set-TextValue "D3" "1.854.80"
Set-TextValue "E3" "  -3.00%  "

Set-TextValue "E4" "  -0.10%  "

Set-TextValue "D5" "329.14"
Set-TextValue "E5" "  +0.42%  "

Set-TextValue "E6" "  -0.06%  "

Set-TextValue "D7" "0.4621"
Set-TextValue "E7" "  -1.09%  "

Set-TextValue "D8" "0.3946"
Set-TextValue "E8" "  -1.49%  "

Set-TextValue "D9" "46.49"
Set-TextValue "E9" "  -12.55%  "

Set-TextValue "D10" "0.07950"
Set-TextValue "E10" "  -5.58%  "

Set-TextValue "E11" "  -3.03%  "

Set-TextValue "D12" "21.49"
Set-TextValue "E12" "  -2.68%  "

Set-TextValue "D13" "1.853.63"
Set-TextValue "E13" "  -3.58%  "

Set-TextValue "D14" "5.923"
Set-TextValue "E14" "  -2.19%  "

Set-TextValue "D15" "7.141"
Set-TextValue "E15" "  -3.70%  "

Set-TextValue "E16" "  -0.03%  "

Set-TextValue "E17" "  -3.95%  "

Set-TextValue "D18" "0.00001032"
Set-TextValue "E18" "  -3.27%  "

Set-TextValue "D19" "0.06578"
Set-TextValue "E19" "  -0.54%  "

Set-TextValue "D20" "17.25"
Set-TextValue "E20" "  -3.79%  "

Set-TextValue "E21" "  -0.08%  "

Set-TextValue "D22" "5.467"
Set-TextValue "E22" "  -4.49%  "

Set-TextValue "D23" "27.324.90"
Set-TextValue "E23" "  -3.04%  "

Set-TextValue "D24" "10.89"
Set-TextValue "E24" "  -2.71%  "

Set-TextValue "D25" "2.304"
Set-TextValue "E25" "  +0.65%  "

Set-TextValue "D26" "2.079.01"
Set-TextValue "E26" "  -3.37%  "

Set-TextValue "B27" "EthereumClassic"
Set-TextValue "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D27" "20.34"
Set-TextValue "E27" "  +1.79%  "

Set-TextValue "B28" "Monero"
Set-TextValue "C28" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D28" "153.78"
Set-TextValue "E28" "  +0.28%  "

Set-TextValue "D29" "2.065"
Set-TextValue "E29" "  -2.50%  "

Set-TextValue "D30" "5.457"
Set-TextValue "E30" "  -4.86%  "

Set-TextValue "D31" "121.75"
Set-TextValue "E31" "  -1.17%  "

Set-TextValue "D32" "0.09438"
Set-TextValue "E32" "  -2.16%  "

Set-TextValue "D33" "0.9494"
Set-TextValue "E33" "  -2.69%  "

Set-TextValue "E34" "  +0.20%  "

Set-TextValue "D35" "3.577"
Set-TextValue "E35" "  -1.92%  "

Set-TextValue "D36" "5.264"
Set-TextValue "E36" "  -4.84%  "

Set-TextValue "D37" "0.06040"
Set-TextValue "E37" "  -2.06%  "

Set-TextValue "D38" "0.02226"
Set-TextValue "E38" "  -2.97%  "

Set-TextValue "E39" "  -4.73%  "

Set-TextValue "E40" "  -0.09%  "

Set-TextValue "D41" "8.019"
Set-TextValue "E41" "  -8.84%  "

Set-TextValue "D42" "0.5928"
Set-TextValue "E42" "  -3.55%  "

Set-TextValue "D43" "0.1890"
Set-TextValue "E43" "  -0.74%  "

Set-TextValue "D44" "10.19"
Set-TextValue "E44" "  -7.47%  "

Set-TextValue "D45" "1.281"
Set-TextValue "E45" "  -2.26%  "

Set-TextValue "D46" "0.5634"
Set-TextValue "E46" "  -3.49%  "

Set-TextValue "D47" "12.13"
Set-TextValue "E47" "  -4.86%  "

Set-TextValue "D48" "3.407"
Set-TextValue "E48" "  -0.58%  "

Set-TextValue "D49" "1.918"
Set-TextValue "E49" "  -4.89%  "

Set-TextValue "D50" "0.06758"
Set-TextValue "E50" "  -1.98%  "

Set-TextValue "D51" "109.29"
Set-TextValue "E51" "  -0.88%  "
